$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the student account usernames and email addresses (rows 2-6):
# "IIT20500xx" / "iit20500xx@iiita.ac.in" -> "IT20500xx" / "it20500xx@thapar.edu"
$ws.Range("B2").Value = "IT2050001"
$ws.Range("C2").Value = "it2050001@thapar.edu"

$ws.Range("B3").Value = "IT2050002"
$ws.Range("C3").Value = "it2050002@thapar.edu"

$ws.Range("B4").Value = "IT2050003"
$ws.Range("C4").Value = "it2050003@thapar.edu"

$ws.Range("B5").Value = "IT2050004"
$ws.Range("C5").Value = "it2050004@thapar.edu"

$ws.Range("B6").Value = "IT2050005"
$ws.Range("C6").Value = "it2050005@thapar.edu"

# Move the active selection, matching the saved workbook state
$ws.Range("C7").Select()

$wb.Save()
